$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "Créditos-trabalho: 4" "Créditos-trabalho: 3"
Replace-Text "Carga horária: 120 h" "Carga horária: 90 h"
Replace-Text "Ativação: 01/01/2016" "Ativação: 01/01/2025"

Replace-Text "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador, o qual deve constituir-se num projeto de engenharia química." "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador,o qual deve constituir-se num projeto de tema específico relacionado às atribuições da profissão."

Replace-Text "Elaboração de uma monografia de conclusão de curso que apresente: (1) o tema e sua importância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) o desenvolvimento do projeto, (6) a análise e discussão dos resultados, (7) as conclusões e (8) referências bibliográficas." "Elaboração de uma monografia ou de relatório técnico que apresente: (1) o tema e suaimportância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) odesenvolvimento do projeto, (6) a análise e discussão dos resultados, (7) as conclusões e (8)referências bibliográficas"

Replace-Text "Reuniões periódicas com o orientador e realização do trabalho de conclusão de curso conforme orientação e apresentação de uma monografia final, conforme norma do Departamento de Engenharia Química." "Reuniões periódicas com o orientador e realização do trabalho conforme orientação eapresentação de uma monografia final, conforme norma do Departamento de Engenharia Químicae Produção."

Replace-Text "Avaliação da monografia perante uma banca examinadora composta por 3 (três) membros, obrigatoriamente docentes da Escola de Engenharia de Lorena (EEL)." "Avaliação perante uma banca examinadora composta por 3 (três) membros, conforme norma doDepartamento de Engenharia Química e Produção."

Replace-Text "Reapresentação da monografia, preferencialmente para a mesma banca, com as modificações sugeridas para uma nova avaliação." "Reavaliação, preferencialmente para a mesma banca, com as modificações sugeridas."
